{"js": "// Apply the \"ubundu notes 2add 18.03.22\" edit:\n//  1. The first (empty) paragraph in the document gets a GitHub token\n//     typed into its existing run.\n//  2. A new, still-empty italic/size-16 \"No Spacing\" paragraph is inserted\n//     right after it (same paragraph-mark formatting as paragraph 1).\n//  3. The \"ls -al // display filenames ...\" line (previously split across\n//     three runs: \"ls -\", \"a\", \"l // display ...\") is normalized back into\n//     a single run with the identical visible text.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// --- 1 & 2: first paragraph gets the token text, then a new blank\n// paragraph is inserted right after it. ---\nconst firstPara = paragraphs.items[0];\nfirstPara.getRange().insertText(\"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\", \"Replace\");\nawait context.sync();\n\n// insertParagraph(\"After\") duplicates the paragraph-mark formatting\n// (italic, 16pt, Times New Roman, \"No Spacing\" style) of firstPara, giving\n// the same empty italic paragraph introduced by the diff.\nfirstPara.insertParagraph(\"\", \"After\");\nawait context.sync();\n\n// --- 3: merge the \"ls -al ...\" runs back into one run. ---\nconst lsAlText =\n  \"ls -al // display filenames with date and time including hidden files also\";\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items/text\");\nawait context.sync();\n\nfor (const p of paragraphs2.items) {\n  if (p.text === lsAlText) {\n    // Re-typing the identical text over the paragraph's own range collapses\n    // the three runs (\"ls -\", \"a\", \"l // display ...\") into a single run\n    // while keeping the same formatting/text.\n    p.getRange().insertText(lsAlText, \"Replace\");\n    break;\n  }\n}\nawait context.sync();\n", "ps1": "# Apply the \"ubundu notes 2add 18.03.22\" edit:\n#  1. The first (empty) paragraph in the document gets a GitHub token\n#     typed into its existing run.\n#  2. A new, still-empty italic/size-16 \"No Spacing\" paragraph is inserted\n#     right after it (same paragraph-mark formatting as paragraph 1).\n#  3. The \"ls -al // display filenames ...\" line (previously split across\n#     three runs: \"ls -\", \"a\", \"l // display ...\") is normalized back into\n#     a single run with the identical visible text.\n\n$d = $word.ActiveDocument\n\n# --- 1 & 2: first paragraph gets the token text, then a new blank\n# paragraph is inserted right after it. ---\n$firstPara = $d.Paragraphs.Item(1)\n$firstPara.Range.InsertBefore(\"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\")\n\n# InsertParagraphAfter duplicates the paragraph-mark formatting (italic,\n# 16pt, Times New Roman, \"No Spacing\" style) of firstPara, giving the same\n# empty italic paragraph introduced by the diff.\n$firstParaAgain = $d.Paragraphs.Item(1)\n$firstParaAgain.Range.InsertParagraphAfter()\n\n# --- 3: merge the \"ls -al ...\" runs back into one run. ---\n$lsAlText = \"ls -al // display filenames with date and time including hidden files also\"\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  $r = $p.Range\n  $r.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark\n  if ($r.Text -eq $lsAlText) {\n    $fontName = $r.Font.Name\n\n    # Deleting and retyping the text collapses the three runs (\"ls -\",\n    # \"a\", \"l // display ...\") into a single run; the font is restored\n    # below since Delete()+InsertBefore() can drop direct formatting.\n    $r.Delete()\n    $r.InsertBefore($lsAlText)\n\n    $r2 = $p.Range\n    $r2.MoveEnd(1, -1) | Out-Null\n    $r2.Font.Name = $fontName\n    $r2.Font.NameBi = $fontName\n    break\n  }\n}\n"}
